$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text content looks like a plain decimal number must be
# forced to Text format *before* assignment, otherwise Excel coerces
# the string into a numeric value (losing trailing zeros / switching to
# scientific notation) instead of keeping the literal text from the feed.
$forceTextValues = [ordered]@{
    'D5' = '577.11'
    'D6' = '136.42'
    'D10' = '7.49'
    'D15' = '0.0000177'
    'D17' = '25.37'
    'D19' = '14.18'
    'D21' = '9.47'
    'D22' = '387.61'
    'D23' = '0.568'
    'D27' = '71.20'
    'D28' = '7.67'
    'D31' = '8.27'
    'D38' = '6.95'
    'D40' = '163.02'
    'D46' = '4.44'
    'D47' = '41.79'
    'D48' = '24.60'
    'D49' = '6.95'
    'D50' = '23.27'
}
foreach ($cell in $forceTextValues.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $forceTextValues[$cell]
}

# Remaining updated cells (plain text, URLs, or strings that already
# fail numeric parsing such as the "thousands-dotted" price strings and
# the padded percentage strings) can be assigned directly.
$plainValues = [ordered]@{
    'D2' = '61.732.13'
    'E2' = '  +2.48%  '
    'D3' = '3.391.19'
    'E3' = '  +1.43%  '
    'E4' = '  +0.01%  '
    'E5' = '  +2.00%  '
    'E6' = '  +3.40%  '
    'E7' = '  +0.09%  '
    'D8' = '3.390.43'
    'E8' = '  +1.46%  '
    'E9' = '  +0.29%  '
    'E10' = '  +0.45%  '
    'E11' = '  +5.78%  '
    'E12' = '  +3.81%  '
    'D13' = '3.972.27'
    'E13' = '  +1.51%  '
    'E14' = '  +2.35%  '
    'E15' = '  +3.93%  '
    'D16' = '3.386.96'
    'E16' = '  +0.74%  '
    'E17' = '  +1.94%  '
    'D18' = '61.763.46'
    'E18' = '  +2.44%  '
    'E19' = '  +4.29%  '
    'E20' = '  +2.26%  '
    'E21' = '  +2.25%  '
    'E23' = '  +1.32%  '
    'D24' = '3.536.84'
    'E24' = '  +1.84%  '
    'E25' = '  -0.04%  '
    'E26' = '  +13.52%  '
    'E27' = '  +2.31%  '
    'E28' = '  +1.90%  '
    'E29' = '  -4.65%  '
    'E30' = '  +0.34%  '
    'E31' = '  +3.31%  '
    'E32' = '  +4.23%  '
    'E33' = '  +1.70%  '
    'E34' = '  +0.04%  '
    'D35' = '3.424.32'
    'E35' = '  +1.52%  '
    'E36' = '  +2.25%  '
    'E37' = '  +0.45%  '
    'E38' = '  +0.77%  '
    'E39' = '  +2.52%  '
    'E40' = '  +2.65%  '
    'E41' = '  +1.52%  '
    'E42' = '  +11.80%  '
    'E43' = '  +4.37%  '
    'E44' = '  +0.02%  '
    'E45' = '  +2.39%  '
    'B46' = 'Filecoin'
    'C46' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'E46' = '  +1.12%  '
    'B47' = 'OKB'
    'C47' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'E47' = '  +2.16%  '
    'E48' = '  +4.58%  '
    'E49' = '  +1.92%  '
    'E50' = '  +3.16%  '
    'D51' = '2.353.82'
    'E51' = '  +7.83%  '
}
foreach ($cell in $plainValues.Keys) {
    $ws.Range($cell).Value = $plainValues[$cell]
}

